$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-08-10 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-08-11 Friday", 2) | Out-Null
$d.Content.Find.Execute("88×42=", $true, $false, $false, $false, $false, $true, 1, $false, "28×31=", 2) | Out-Null
$d.Content.Find.Execute("13×39=", $true, $false, $false, $false, $false, $true, 1, $false, "52×43=", 2) | Out-Null
$d.Content.Find.Execute("57×28=", $true, $false, $false, $false, $false, $true, 1, $false, "12×89=", 2) | Out-Null
$d.Content.Find.Execute("52×58=", $true, $false, $false, $false, $false, $true, 1, $false, "48×41=", 2) | Out-Null
$d.Content.Find.Execute("43×90=", $true, $false, $false, $false, $false, $true, 1, $false, "15×70=", 2) | Out-Null
$d.Content.Find.Execute("26×73=", $true, $false, $false, $false, $false, $true, 1, $false, "54×95=", 2) | Out-Null
$d.Content.Find.Execute("62×59=", $true, $false, $false, $false, $false, $true, 1, $false, "39×46=", 2) | Out-Null
$d.Content.Find.Execute("97×71=", $true, $false, $false, $false, $false, $true, 1, $false, "85×36=", 2) | Out-Null
$d.Content.Find.Execute("40×91=", $true, $false, $false, $false, $false, $true, 1, $false, "83×89=", 2) | Out-Null
$d.Content.Find.Execute("38×45=", $true, $false, $false, $false, $false, $true, 1, $false, "13×27=", 2) | Out-Null
$d.Content.Find.Execute("98×62=", $true, $false, $false, $false, $false, $true, 1, $false, "86×42=", 2) | Out-Null
$d.Content.Find.Execute("77×45=", $true, $false, $false, $false, $false, $true, 1, $false, "63×76=", 2) | Out-Null
$d.Content.Find.Execute("84×86=", $true, $false, $false, $false, $false, $true, 1, $false, "43×82=", 2) | Out-Null
$d.Content.Find.Execute("42×28=", $true, $false, $false, $false, $false, $true, 1, $false, "81×89=", 2) | Out-Null
$d.Content.Find.Execute("47×29=", $true, $false, $false, $false, $false, $true, 1, $false, "40×54=", 2) | Out-Null
$d.Content.Find.Execute("80×47=", $true, $false, $false, $false, $false, $true, 1, $false, "14×44=", 2) | Out-Null
$d.Content.Find.Execute("29×37=", $true, $false, $false, $false, $false, $true, 1, $false, "30×40=", 2) | Out-Null
$d.Content.Find.Execute("65×58=", $true, $false, $false, $false, $false, $true, 1, $false, "15×95=", 2) | Out-Null
$d.Content.Find.Execute("99×80=", $true, $false, $false, $false, $false, $true, 1, $false, "72×97=", 2) | Out-Null
$d.Content.Find.Execute("43×86=", $true, $false, $false, $false, $false, $true, 1, $false, "39×57=", 2) | Out-Null
$d.Content.Find.Execute("85×28=", $true, $false, $false, $false, $false, $true, 1, $false, "42×52=", 2) | Out-Null
$d.Content.Find.Execute("25×91=", $true, $false, $false, $false, $false, $true, 1, $false, "38×88=", 2) | Out-Null
$d.Content.Find.Execute("91×32=", $true, $false, $false, $false, $false, $true, 1, $false, "21×69=", 2) | Out-Null
$d.Content.Find.Execute("77×53=", $true, $false, $false, $false, $false, $true, 1, $false, "11×38=", 2) | Out-Null
$d.Content.Find.Execute("19×75=", $true, $false, $false, $false, $false, $true, 1, $false, "39×89=", 2) | Out-Null
